# Weekly update: two new price records (most recent week's data) are
# inserted into the daily log, pushing the older rows down by one each
# time. Net result: 20 data rows (rows 2-21) become 22 data rows
# (rows 2-23); dimension grows from A1:R21 to A1:R23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new row at row 5 (shifts old rows 5-21 down to 6-22)
$ws.Rows("5:5").Insert()

# Insert the second new row at row 13 (shifts the now-current rows
# 13-22 down to 14-23)
$ws.Rows("13:13").Insert()

# --- Fill in the brand-new row 5 ---
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value = "Metropolitana"
$ws.Cells.Item(5, 4).Value = 44575
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = 100112010
$ws.Cells.Item(5, 7).Value = "Achicoria"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 61
$ws.Cells.Item(5, 11).Value = 8000
$ws.Cells.Item(5, 12).Value = 8000
$ws.Cells.Item(5, 13).Value = 8000
$ws.Cells.Item(5, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(5, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(5, 16).Value = 500
$ws.Cells.Item(5, 17).Value = 16
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# --- Fill in the brand-new row 13 ---
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44573
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = 100112010
$ws.Cells.Item(13, 7).Value = "Achicoria"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 34
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 8000
$ws.Cells.Item(13, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(13, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 16).Value = 500
$ws.Cells.Item(13, 17).Value = 16
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Make sure the date cells carry the same date/time number format as
# the rest of column D.
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("D13").NumberFormat = $ws.Range("D14").NumberFormat

$ws.Range("A1").Select()
